# Update "paises" (countries) COVID tracking sheet:
#  - bump the "Datos actualizados" timestamp
#  - refresh per-country statistics for a new data pull
#  - re-sort: "Emiratos Arabes Unidos" overtakes "Paises Bajos" and
#    "Bielorrusia" in total cases, so those three rows swap places
#    (data for Paises Bajos/Bielorrusia shifts down a row, and the
#    now-top row gets Emiratos Arabes Unidos' refreshed numbers)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Septiembre de 2020 a las 12:33"

# Row 4 - Estados Unidos
$ws.Cells.Item(4, 2).Value = 6335653
$ws.Cells.Item(4, 3).Value = 409
$ws.Cells.Item(4, 4).Value = 3575495
$ws.Cells.Item(4, 5).Value = 2569098
$ws.Cells.Item(4, 7).Value = 2
$ws.Cells.Item(4, 8).Value = 191060

# Row 6 - India
$ws.Cells.Item(6, 2).Value = 3940131
$ws.Cells.Item(6, 3).Value = 7007
$ws.Cells.Item(6, 5).Value = 834382
$ws.Cells.Item(6, 7).Value = 29
$ws.Cells.Item(6, 8).Value = 68598

# Row 15 - Iran
$ws.Cells.Item(15, 2).Value = 382772
$ws.Cells.Item(15, 3).Value = 2026
$ws.Cells.Item(15, 4).Value = 330308
$ws.Cells.Item(15, 5).Value = 30420
$ws.Cells.Item(15, 7).Value = 118
$ws.Cells.Item(15, 8).Value = 22044

# Row 17 - Banglades
$ws.Cells.Item(17, 2).Value = 321615
$ws.Cells.Item(17, 3).Value = 1929
$ws.Cells.Item(17, 4).Value = 216191
$ws.Cells.Item(17, 5).Value = 101012
$ws.Cells.Item(17, 7).Value = 29
$ws.Cells.Item(17, 8).Value = 4412

# Row 29 - Israel
$ws.Cells.Item(29, 2).Value = 125755
$ws.Cells.Item(29, 3).Value = 1300
$ws.Cells.Item(29, 4).Value = 99487
$ws.Cells.Item(29, 5).Value = 25277
$ws.Cells.Item(29, 7).Value = 6
$ws.Cells.Item(29, 8).Value = 991

# Row 37 - Rumania
$ws.Cells.Item(37, 2).Value = 92595
$ws.Cells.Item(37, 3).Value = 1339
$ws.Cells.Item(37, 4).Value = 39626
$ws.Cells.Item(37, 5).Value = 49157
$ws.Cells.Item(37, 7).Value = 47
$ws.Cells.Item(37, 8).Value = 3812

# Rows 44-46 re-sort: Emiratos Arabes Unidos moves above Paises Bajos and
# Bielorrusia in the ranking. Update the country labels (column A) to match
# the new order, and move each row's statistics along with its country.
# New row 44 (Emiratos Arabes Unidos) gets freshly updated numbers.
$ws.Cells.Item(44, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(44, 2).Value = 72766
$ws.Cells.Item(44, 3).Value = 612
$ws.Cells.Item(44, 4).Value = 63158
$ws.Cells.Item(44, 5).Value = 9221
$ws.Cells.Item(44, 8).Value = 387

# New row 45 (Paises Bajos) takes what used to be row 44's values.
$ws.Cells.Item(45, 1).Value = "Paises Bajos"
$ws.Cells.Item(45, 2).Value = 72464
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(45, 5).Value = 0
$ws.Cells.Item(45, 8).Value = 6235

# New row 46 (Bielorrusia) takes what used to be row 45's values.
$ws.Cells.Item(46, 1).Value = "Bielorrusia"
$ws.Cells.Item(46, 2).Value = 72302
$ws.Cells.Item(46, 4).Value = 71205
$ws.Cells.Item(46, 5).Value = 401
$ws.Cells.Item(46, 8).Value = 696

# Row 55 - Barein
$ws.Cells.Item(55, 5).Value = 3228
$ws.Cells.Item(55, 7).Value = 2
$ws.Cells.Item(55, 8).Value = 192

# Row 62 - Suiza
$ws.Cells.Item(62, 2).Value = 43532
$ws.Cells.Item(62, 3).Value = 405
$ws.Cells.Item(62, 5).Value = 5019

# Row 71 - Austria
$ws.Cells.Item(71, 2).Value = 28729
$ws.Cells.Item(71, 3).Value = 357
$ws.Cells.Item(71, 4).Value = 24513
$ws.Cells.Item(71, 5).Value = 3481

# Row 87 - Senegal
$ws.Cells.Item(87, 2).Value = 13881
$ws.Cells.Item(87, 3).Value = 55
$ws.Cells.Item(87, 4).Value = 9723
$ws.Cells.Item(87, 5).Value = 3869
$ws.Cells.Item(87, 7).Value = 2
$ws.Cells.Item(87, 8).Value = 289

# Row 102 - Finlandia
$ws.Cells.Item(102, 2).Value = 8225
$ws.Cells.Item(102, 3).Value = 25
$ws.Cells.Item(102, 5).Value = 539

# Row 130 - Eslovenia
$ws.Cells.Item(130, 2).Value = 3079
$ws.Cells.Item(130, 3).Value = 47
$ws.Cells.Item(130, 4).Value = 2440
$ws.Cells.Item(130, 5).Value = 505

# Row 149 - Nueva Zelanda
$ws.Cells.Item(149, 5).Value = 111
$ws.Cells.Item(149, 7).Value = 1
$ws.Cells.Item(149, 8).Value = 23

# Row 184 - Gibraltar
$ws.Cells.Item(184, 2).Value = 305
$ws.Cells.Item(184, 3).Value = 7
$ws.Cells.Item(184, 4).Value = 248
$ws.Cells.Item(184, 5).Value = 57
